$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*admin dog detailsnál ha dog status adopted, örökbe adás gomb legyen disabled*") {
        $p.Range.Delete()
        break
    }
}
